$wb = $excel.ActiveWorkbook

# --- Sheet 2 ("Kế Hoạch") content updates ---
$wsPlan = $wb.Worksheets.Item("Kế Hoạch")

# Row 23: person responsible changes from "Cả team" to "Phong + Trường"
$wsPlan.Range("C23").Value = "Phong + Trường"

# Row 24 (new): continuation of "Thiết kế web" task -> "Code chức năng"
$wsPlan.Range("B24").Value = "Code chức năng"
$wsPlan.Range("C24").Value = "Hoàng"

# Row 25 (new): week 5 - "Thiết kế web" / "Phong + Trường"
$wsPlan.Range("A25").Value = 5
$wsPlan.Range("A25").HorizontalAlignment = -4108
$wsPlan.Range("B25").Value = "Thiết kế web"
$wsPlan.Range("C25").Value = "Phong + Trường"

# Row 26 (new): continuation -> "Code chức năng"
$wsPlan.Range("B26").Value = "Code chức năng"
$wsPlan.Range("C26").Value = "Hoàng"

# Row 27 (new): week 6 - "Hoàn thiện web" / "Cả team (Hoàng)"
$wsPlan.Range("A27").Value = 6
$wsPlan.Range("A27").HorizontalAlignment = -4108
$wsPlan.Range("B27").Value = "Hoàn thiện web"
$wsPlan.Range("C27").Value = "Cả team (Hoàng)"

# Row 28 (new): week 7 - "Hoàn thiện web" / "Cả team (Hoàng)"
$wsPlan.Range("A28").Value = 7
$wsPlan.Range("A28").HorizontalAlignment = -4108
$wsPlan.Range("B28").Value = "Hoàn thiện web"
$wsPlan.Range("C28").Value = "Cả team (Hoàng)"

# --- View state ---
# "Danh sách Nhóm + Tuần" keeps a scrolled-down view with a different selection
$wsGroup = $wb.Worksheets.Item("Danh sách Nhóm + Tuần")
$wsGroup.Activate()
$wsGroup.Range("I23").Select()

# "Kế Hoạch" becomes the active / selected sheet (done last so it ends up active)
$wsPlan.Activate()
$wsPlan.Range("B29").Select()
